# Updated cryptos list on Sat Jul 13 04:49:53 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '57.794.15'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +1.44%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.114.27'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.10%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '532.44'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.38%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.16'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +1.55%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.25%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.498'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +10.47%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.34'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.08%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.34%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.412'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.140'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +3.54%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.652.18'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +1.22%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.62'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000167'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +3.72%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '57.886.43'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.50%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.129.12'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +1.80%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.13'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +4.26%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '12.79'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +2.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.10'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +3.18%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '374.16'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +7.68%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -2.23%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.41'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +1.65%  '
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +1.87%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.999'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0₃0879'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +1.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.62'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +4.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.12'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +4.19%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.49'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +3.40%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.13'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.21%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +2.95%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '160.49'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.90%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.17'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +3.03%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.20%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '25.44'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +3.79%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0669'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +2.76%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.563.10'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +7.49%  '
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  +3.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '38.55'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +5.36%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.696'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.01%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +2.43%  '
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.977'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +2.33%  '
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +3.40%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '19.88'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0953'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +7.05%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.746'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.64%  '
